# Apply a row-permutation of runs/balls/fours/sixes data (columns C:F, rows 2-17)
# on the "David Warner (c)" sheet. Values are stored as text (numberStoredAsText),
# so we assign them as strings to preserve that representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (r => runs, balls, fours, sixes), matching the target diff.
$data = @{
    2  = @("2", "3", "0", "0")
    3  = @("36", "30", "2", "1")
    4  = @("85", "58", "10", "1")
    5  = @("47", "33", "5", "0")
    6  = @("52", "40", "5", "1")
    7  = @("45", "33", "3", "2")
    8  = @("28", "29", "3", "0")
    9  = @("4", "4", "1", "0")
    10 = @("9", "13", "0", "0")
    11 = @("48", "38", "3", "2")
    12 = @("17", "17", "3", "0")
    13 = @("8", "5", "0", "1")
    14 = @("66", "34", "8", "2")
    15 = @("35", "20", "3", "2")
    16 = @("6", "6", "1", "0")
    17 = @("60", "44", "5", "2")
}

# Ensure the target range keeps a Text number format so values remain
# stored as strings (matching the workbook's original numberStoredAsText cells).
$ws.Range("C2:F17").NumberFormat = "@"

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
    $ws.Cells.Item($r, 6).Value = $vals[3]
}
